$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the cryptos table with the latest scraped prices/volumes.
# (Generated from the per-row OOXML diff; D-column price cells are forced
#  to text format so values like "63.103.56" or "416.23" are preserved
#  exactly instead of being auto-parsed as numbers by Excel.)

# Row 2
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "63.103.56"
$ws.Range("E2").Value = "  +10.24%  "

# Row 3
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.474.07"
$ws.Range("E3").Value = "  +6.54%  "

# Row 4
$ws.Range("E4").Value = "  +0.02%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "416.23"
$ws.Range("E5").Value = "  +4.59%  "

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "123.00"
$ws.Range("E6").Value = "  +12.90%  "

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "3.467.78"
$ws.Range("E7").Value = "  +6.53%  "

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.598"
$ws.Range("E8").Value = "  +3.23%  "

# Row 9
$ws.Range("E9").Value = "  +0.05%  "

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.696"
$ws.Range("E10").Value = "  +12.50%  "

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.132"
$ws.Range("E11").Value = "  +38.60%  "

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "41.53"
$ws.Range("E12").Value = "  +5.91%  "

# Row 13
$ws.Range("E13").Value = "  +0.40%  "

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "4.024.50"
$ws.Range("E14").Value = "  +6.59%  "

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "8.65"
$ws.Range("E15").Value = "  +4.98%  "

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "20.11"
$ws.Range("E16").Value = "  +6.18%  "

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "3.474.72"
$ws.Range("E17").Value = "  +6.36%  "

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "62.924.21"
$ws.Range("E18").Value = "  +10.27%  "

# Row 19
$ws.Range("E19").Value = "  +0.49%  "

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "10.96"
$ws.Range("E20").Value = "  -1.14%  "

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.0000140"
$ws.Range("E21").Value = "  +30.64%  "

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "3.32"
$ws.Range("E22").Value = "  +0.44%  "

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "318.43"
$ws.Range("E23").Value = "  +6.94%  "

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "82.18"
$ws.Range("E24").Value = "  +11.08%  "

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "12.92"
$ws.Range("E25").Value = "  +0.16%  "

# Row 26
$ws.Range("E26").Value = "  +0.29%  "

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "31.13"
$ws.Range("E27").Value = "  +10.76%  "

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "7.80"
$ws.Range("E28").Value = "  +4.84%  "

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "7.90"
$ws.Range("E29").Value = "  +0.63%  "

# Row 30
$ws.Range("B30").Value = "Kaspa"
$ws.Range("C30").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.175"
$ws.Range("E30").Value = "  +4.03%  "

# Row 31
$ws.Range("B31").Value = "LEO"
$ws.Range("C31").Value = "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "4.30"
$ws.Range("E31").Value = "  -1.74%  "

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.117"
$ws.Range("E32").Value = "  +4.27%  "

# Row 33
$ws.Range("E33").Value = "  +4.11%  "

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "42.21"
$ws.Range("E34").Value = "  +5.77%  "

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "2.56"
$ws.Range("E35").Value = "  +19.61%  "

# Row 36
$ws.Range("E36").Value = "  +0.63%  "

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.0491"
$ws.Range("E37").Value = "  -0.77%  "

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "52.29"
$ws.Range("E38").Value = "  +1.83%  "

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.997"
$ws.Range("E39").Value = "  -0.21%  "

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "3.50"
$ws.Range("E40").Value = "  +1.15%  "

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "3.04"
$ws.Range("E41").Value = "  +0.36%  "

# Row 42
$ws.Range("E42").Value = "  +7.10%  "

# Row 43
$ws.Range("E43").Value = "  +3.86%  "

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "136.00"
$ws.Range("E44").Value = "  -1.32%  "

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "16.94"
$ws.Range("E45").Value = "  +1.43%  "

# Row 46
$ws.Range("E46").Value = "  -1.27%  "

# Row 47
$ws.Range("E47").Value = "  +0.70%  "

# Row 48
$ws.Range("E48").Value = "  +2.56%  "

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "22.06"
$ws.Range("E49").Value = "  -1.02%  "

# Row 50
$ws.Range("B50").Value = "Maker"
$ws.Range("C50").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "2.185.52"
$ws.Range("E50").Value = "  +1.89%  "

# Row 51
$ws.Range("B51").Value = "ApeXProtocol"
$ws.Range("C51").Value = "https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "2.47"
$ws.Range("E51").Value = "  +0.31%  "
